$wb = $excel.ActiveWorkbook
for ($i=1; $i -le $wb.Worksheets.Count; $i++) {
    $s = $wb.Worksheets.Item($i)
    try {
        $s.Name = "!" + $s.Name
    } catch {
        Write-Host "Failed on $i : $_"
    }
}
